$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2102.0435
$ws.Range("I15").Value = 2102.0435
$ws.Range("K15").Value = 6306.130500000001
$ws.Range("M15").Value = -6137.130500000001

$ws.Range("H32").Value = 794.6667
$ws.Range("I32").Value = 1250
$ws.Range("J32").Value = 567
$ws.Range("K32").Value = 1250
$ws.Range("L32").Value = 567
$ws.Range("M32").Value = -924
$ws.Range("N32").Value = -1219

$ws.Range("H33").Value = 767.8
$ws.Range("I33").Value = 321.8889
$ws.Range("J33").Value = 1436.6666
$ws.Range("K33").Value = 321.8889
$ws.Range("L33").Value = 1436.6666
$ws.Range("M33").Value = -92.88889999999998
$ws.Range("N33").Value = -1894.6666

$ws.Range("H86").Value = 45476424
$ws.Range("I86").Value = 12000
$ws.Range("J86").Value = 47641396
$ws.Range("K86").Value = 12000
$ws.Range("L86").Value = 47641396
$ws.Range("M86").Value = -10877
$ws.Range("N86").Value = -47643642

$ws.Range("H89").Value = 45476424
$ws.Range("I89").Value = 12000
$ws.Range("J89").Value = 47641396
$ws.Range("K89").Value = 60000
$ws.Range("L89").Value = 238206980
$ws.Range("M89").Value = -54384
$ws.Range("N89").Value = -238218212

$ws.Range("H98").Value = 32039.625
$ws.Range("I98").Value = 51269.727
$ws.Range("J98").Value = 15768
$ws.Range("K98").Value = 51269.727
$ws.Range("L98").Value = 15768
$ws.Range("M98").Value = -49771.727
$ws.Range("N98").Value = -18764

$ws.Range("H105").Value = 37500
$ws.Range("J105").Value = 37500
$ws.Range("L105").Value = 37500
$ws.Range("N105").Value = -44488

$ws.Range("H122").Value = 32039.625
$ws.Range("I122").Value = 51269.727
$ws.Range("J122").Value = 15768
$ws.Range("K122").Value = 153809.181
$ws.Range("L122").Value = 47304
$ws.Range("M122").Value = -151359.181
$ws.Range("N122").Value = -52204

$ws.Range("H137").Value = 12159.1
$ws.Range("I137").Value = 16906.691
$ws.Range("J137").Value = 3342.1428
$ws.Range("K137").Value = 50720.073
$ws.Range("L137").Value = 10026.4284
$ws.Range("M137").Value = -48170.073
$ws.Range("N137").Value = -15126.4284

$ws.Range("H138").Value = 2422.0881
$ws.Range("J138").Value = 3836.606
$ws.Range("L138").Value = 11509.818
$ws.Range("N138").Value = -21789.818

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 86.666664
$ws.Range("I4").Value = 86.666664
$ws.Range("K4").Value = 86.666664
$ws.Range("M4").Value = 29.333336

$ws.Range("H5").Value = 681.8125
$ws.Range("I5").Value = 694.26666
$ws.Range("K5").Value = 694.26666
$ws.Range("M5").Value = -582.26666

$ws.Range("H45").Value = 78638.03999999999
$ws.Range("I45").Value = 147452.22
$ws.Range("J45").Value = 4530.4614
$ws.Range("K45").Value = 147452.22
$ws.Range("L45").Value = 4530.4614
$ws.Range("M45").Value = -147075.22
$ws.Range("N45").Value = -5284.4614

$ws.Range("H96").Value = 27500
$ws.Range("J96").Value = 27500
$ws.Range("L96").Value = 27500
$ws.Range("N96").Value = -32992

$ws.Range("H122").Value = 1697134.5
$ws.Range("I122").Value = 4936.6875
$ws.Range("J122").Value = 4404651
$ws.Range("K122").Value = 14810.0625
$ws.Range("L122").Value = 13213953
$ws.Range("M122").Value = -12360.0625
$ws.Range("N122").Value = -13218853

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 681.8125
$ws.Range("I4").Value = 694.26666
$ws.Range("K4").Value = 694.26666
$ws.Range("M4").Value = -579.26666

$ws.Range("H20").Value = 3444.7368
$ws.Range("I20").Value = 2613.125
$ws.Range("J20").Value = 4049.5454
$ws.Range("K20").Value = 2613.125
$ws.Range("L20").Value = 4049.5454
$ws.Range("M20").Value = -2366.125
$ws.Range("N20").Value = -4543.5454

$ws.Range("H134").Value = 8207.789000000001
$ws.Range("I134").Value = 8910.352999999999
$ws.Range("K134").Value = 26731.059
$ws.Range("M134").Value = -24196.059

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6015.647
$ws.Range("I31").Value = 5914.8965
$ws.Range("K31").Value = 5914.8965
$ws.Range("M31").Value = -5619.8965

$ws.Range("H34").Value = 6015.647
$ws.Range("I34").Value = 5914.8965
$ws.Range("K34").Value = 5914.8965
$ws.Range("M34").Value = -5712.8965

$ws.Range("H58").Value = 3466.2273
$ws.Range("I58").Value = 3995.818
$ws.Range("J58").Value = 2936.6365
$ws.Range("K58").Value = 3995.818
$ws.Range("L58").Value = 2936.6365
$ws.Range("M58").Value = -3792.818
$ws.Range("N58").Value = -3342.6365

$ws.Range("H122").Value = 9053.6875
$ws.Range("J122").Value = 1865
$ws.Range("L122").Value = 5595
$ws.Range("N122").Value = -10495

$ws.Range("H136").Value = 3466.2273
$ws.Range("I136").Value = 3995.818
$ws.Range("J136").Value = 2936.6365
$ws.Range("K136").Value = 11987.454
$ws.Range("L136").Value = 8809.9095
$ws.Range("M136").Value = -9437.454000000002
$ws.Range("N136").Value = -13909.9095

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 256.3
$ws.Range("I8").Value = 256.3
$ws.Range("K8").Value = 768.9000000000001
$ws.Range("M8").Value = -629.9000000000001

$ws.Range("H14").Value = 423.66666
$ws.Range("I14").Value = 423.66666
$ws.Range("K14").Value = 1270.99998
$ws.Range("M14").Value = -1097.99998

$ws.Range("H107").Value = 1040.4231
$ws.Range("I107").Value = 358.33334
$ws.Range("J107").Value = 1401.5294
$ws.Range("K107").Value = 1075.00002
$ws.Range("L107").Value = 4204.5882
$ws.Range("M107").Value = 844.9999800000001
$ws.Range("N107").Value = -8044.5882

$ws.Range("H131").Value = 1607.6224
$ws.Range("I131").Value = 1732.8334
$ws.Range("K131").Value = 5198.5002
$ws.Range("M131").Value = -158.5002000000004

$ws.Range("H137").Value = 2075.889
$ws.Range("I137").Value = 2155.0625
$ws.Range("J137").Value = 1442.5
$ws.Range("K137").Value = 6465.1875
$ws.Range("L137").Value = 4327.5
$ws.Range("M137").Value = -1365.1875
$ws.Range("N137").Value = -14527.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 119.92308
$ws.Range("I2").Value = 85.09090999999999
$ws.Range("J2").Value = 311.5
$ws.Range("K2").Value = 85.09090999999999
$ws.Range("L2").Value = 311.5
$ws.Range("M2").Value = 27.90909000000001
$ws.Range("N2").Value = -537.5

$ws.Range("H45").Value = 32199.75
$ws.Range("J45").Value = 32199.75
$ws.Range("L45").Value = 32199.75
$ws.Range("N45").Value = -33317.75

$ws.Range("H70").Value = 5948.6924
$ws.Range("I70").Value = 5417.6816
$ws.Range("K70").Value = 5417.6816
$ws.Range("M70").Value = -5147.6816

$ws.Range("H73").Value = 5948.6924
$ws.Range("I73").Value = 5417.6816
$ws.Range("K73").Value = 5417.6816
$ws.Range("M73").Value = -4481.6816

$ws.Range("H113").Value = 9772.467000000001
$ws.Range("I113").Value = 17321.428
$ws.Range("J113").Value = 3167.125
$ws.Range("K113").Value = 17321.428
$ws.Range("L113").Value = 3167.125
$ws.Range("M113").Value = -15151.428
$ws.Range("N113").Value = -7507.125

$ws.Range("H122").Value = 9058.929
$ws.Range("I122").Value = 6547.6665
$ws.Range("J122").Value = 13579.2
$ws.Range("K122").Value = 19642.9995
$ws.Range("L122").Value = 40737.60000000001
$ws.Range("M122").Value = -17192.9995
$ws.Range("N122").Value = -45637.60000000001

$ws.Range("H132").Value = 2242.426
$ws.Range("I132").Value = 2173.35
$ws.Range("K132").Value = 6520.049999999999
$ws.Range("M132").Value = -3990.049999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 26180.857
$ws.Range("I7").Value = 35485.785
$ws.Range("K7").Value = 35485.785
$ws.Range("M7").Value = -35373.785

$ws.Range("H126").Value = 26180.857
$ws.Range("I126").Value = 35485.785
$ws.Range("K126").Value = 106457.355
$ws.Range("M126").Value = -103987.355

$ws.Range("H135").Value = 94826.2
$ws.Range("J135").Value = 94826.2
$ws.Range("L135").Value = 94826.2
$ws.Range("N135").Value = -104966.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3932.7068
$ws.Range("I122").Value = 1599.8235
$ws.Range("J122").Value = 7237.625
$ws.Range("K122").Value = 4799.470499999999
$ws.Range("L122").Value = 21712.875
$ws.Range("M122").Value = -2349.470499999999
$ws.Range("N122").Value = -26612.875

$ws.Range("H124").Value = 45429
$ws.Range("J124").Value = 45429
$ws.Range("L124").Value = 45429
$ws.Range("N124").Value = -55249

$ws.Range("H132").Value = 9195.714
$ws.Range("J132").Value = 4381.467
$ws.Range("L132").Value = 13144.401
$ws.Range("N132").Value = -18204.401
